$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("J6").Value = 1.73
$ws.Range("K6").Value = 2.5
$ws.Range("Q6").Value = 1.77
$ws.Range("R6").Value = 2
$ws.Range("U6").Value = 2.25
$ws.Range("V6").Value = 1.57
$ws.Range("W6").Value = 6.5
$ws.Range("AC6").Value = 11
$ws.Range("AW6").Value = 11
$ws.Range("Q7").Value = 1.69
$ws.Range("R7").Value = 2.07
$ws.Range("M8").Value = 1.05
$ws.Range("N8").Value = 11
$ws.Range("Q9").Value = 2.15
$ws.Range("R9").Value = 1.61
$ws.Range("G13").Value = 3.6
$ws.Range("H13").Value = 3
$ws.Range("J13").Value = 4.5
$ws.Range("K13").Value = 1.95
$ws.Range("M13").Value = 1.1
$ws.Range("N13").Value = 7
$ws.Range("O13").Value = 1.5
$ws.Range("P13").Value = 2.5
$ws.Range("Q13").Value = 2.6
$ws.Range("R13").Value = 1.48
$ws.Range("S13").Value = 1.54
$ws.Range("T13").Value = 2.25
$ws.Range("U13").Value = 2.1
$ws.Range("V13").Value = 1.67
$ws.Range("AC13").Value = 6.5
$ws.Range("AE13").Value = 19
$ws.Range("AK13").Value = 21
$ws.Range("AL13").Value = 41
$ws.Range("AM13").Value = 501
$ws.Range("AO13").Value = 23
$ws.Range("AR13").Value = 126
$ws.Range("AS13").Value = 351
$ws.Range("AT13").Value = 2.25
$ws.Range("AX13").Value = 13
$ws.Range("AY13").Value = 29
$ws.Range("BA13").Value = 81
$ws.Range("BD13").Value = 151
$ws.Range("S14").Value = 1.33
$ws.Range("J15").Value = 2.63
$ws.Range("K15").Value = 2.1
$ws.Range("N15").Value = 9.5
$ws.Range("S15").Value = 1.37
$ws.Range("U15").Value = 1.83
$ws.Range("V15").Value = 1.83
$ws.Range("W15").Value = 7
$ws.Range("AA15").Value = 17
$ws.Range("AB15").Value = 29
$ws.Range("AC15").Value = 9.5
$ws.Range("G16").Value = 2.4
$ws.Range("H16").Value = 3.8
$ws.Range("I16").Value = 2.5
$ws.Range("J16").Value = 2.77
$ws.Range("K16").Value = 2.45
$ws.Range("L16").Value = 2.9
$ws.Range("M16").Value = 1.02
$ws.Range("S16").Value = 1.23
$ws.Range("T16").Value = 3.75
$ws.Range("V16").Value = 2.82
$ws.Range("W16").Value = 15
$ws.Range("X16").Value = 17.5
$ws.Range("Z16").Value = 29
$ws.Range("AD16").Value = 8.5
$ws.Range("AE16").Value = 10.75
$ws.Range("AF16").Value = 29
$ws.Range("AG16").Value = 15
$ws.Range("AH16").Value = 17.5
$ws.Range("AI16").Value = 10
$ws.Range("AJ16").Value = 30
$ws.Range("AP16").Value = 14
$ws.Range("AR16").Value = 45
$ws.Range("AT16").Value = 3.75
$ws.Range("AW16").Value = 5.1
$ws.Range("AY16").Value = 15
$ws.Range("BB16").Value = 120
